# "adicionado idade ao sistema de login"
#
# Append, after the existing final paragraph ("sistema de login"):
#   1. a blank paragraph
#   2. a paragraph containing "adicione sua idade"
# Both new paragraphs inherit the same paragraph/character formatting
# (center alignment, Times New Roman, size 22) that the final paragraph
# of the document already uses.

$d = $word.ActiveDocument

# Start from the end of the document's last paragraph ("sistema de login").
$lastPara = $d.Paragraphs.Last
$tail = $lastPara.Range
$tail.Collapse(0)

# Insert the first new (blank) paragraph.
$tail.InsertParagraphAfter()
$tail.Collapse(0)

# Insert the second new paragraph and fill it with the requested text.
$tail.InsertParagraphAfter()
$tail.Collapse(0)

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "adicione sua idade"
